$wb = $excel.ActiveWorkbook

# Sheet1
$ws = $wb.Worksheets.Item(1)
$ws.Range("H12").Value = 136.75
$ws.Range("J12").Value = 148.75
$ws.Range("L12").Value = 148.75
$ws.Range("N12").Value = -488.75
$ws.Range("H15").Value = 3799.5356
$ws.Range("I15").Value = 3799.5356
$ws.Range("K15").Value = 11398.6068
$ws.Range("M15").Value = -11229.6068
$ws.Range("H28").Value = 1325.3334
$ws.Range("I28").Value = 1389.2
$ws.Range("K28").Value = 1389.2
$ws.Range("M28").Value = -904.2
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").Value = ""
$ws.Range("H32").Value = 39998.332
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").Value = ""
$ws.Range("H41").Value = 1124.1177
$ws.Range("I41").Value = 954.125
$ws.Range("J41").Value = 1275.2222
$ws.Range("K41").Value = 954.125
$ws.Range("L41").Value = 1275.2222
$ws.Range("M41").Value = -514.125
$ws.Range("N41").Value = -2155.2222
$ws.Range("H53").Value = 410.94736
$ws.Range("I53").Value = 246.45454
$ws.Range("J53").Value = 637.125
$ws.Range("K53").Value = 246.45454
$ws.Range("L53").Value = 637.125
$ws.Range("M53").Value = 390.54546
$ws.Range("N53").Value = -1911.125
$ws.Range("H57").Value = 29144.75
$ws.Range("J57").Value = 29144.75
$ws.Range("L57").Value = 87434.25
$ws.Range("N57").Value = -88432.25
$ws.Range("H62").Value = 3585.3333
$ws.Range("I62").Value = 3377.5
$ws.Range("K62").Value = 3377.5
$ws.Range("M62").Value = -2753.5
$ws.Range("H65").Value = 3585.3333
$ws.Range("I65").Value = 3377.5
$ws.Range("K65").Value = 16887.5
$ws.Range("M65").Value = -13767.5
$ws.Range("H69").Value = 70015
$ws.Range("J69").Value = 70015
$ws.Range("L69").Value = 210045
$ws.Range("N69").Value = -211793
$ws.Range("H72").Value = 70015
$ws.Range("J72").Value = 70015
$ws.Range("L72").Value = 630135
$ws.Range("N72").Value = -638871
$ws.Range("H74").Value = 29974.25
$ws.Range("I74").Value = 29974.25
$ws.Range("K74").Value = 29974.25
$ws.Range("M74").Value = -29038.25
$ws.Range("H77").Value = 29974.25
$ws.Range("I77").Value = 29974.25
$ws.Range("K77").Value = 149871.25
$ws.Range("M77").Value = -145191.25
$ws.Range("H132").Value = 2564.7878
$ws.Range("I132").Value = 2287.9666
$ws.Range("K132").Value = 6863.899800000001
$ws.Range("M132").Value = -4333.899800000001
$ws.Range("H137").Value = 3250.8667
$ws.Range("I137").Value = 1411.375
$ws.Range("J137").Value = 5353.143
$ws.Range("K137").Value = 4234.125
$ws.Range("L137").Value = 16059.429
$ws.Range("M137").Value = -1684.125
$ws.Range("N137").Value = -21159.429
$ws.Range("H138").Value = 2526.0293
$ws.Range("I138").Value = 896.1
$ws.Range("J138").Value = 3205.1667
$ws.Range("K138").Value = 2688.3
$ws.Range("L138").Value = 9615.500100000001
$ws.Range("M138").Value = 2451.7
$ws.Range("N138").Value = -19895.5001
$ws.Range("H139").Value = 48000
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = ""
$ws.Range("H141").Value = 5603.8823
$ws.Range("I141").Value = 5417.8
$ws.Range("K141").Value = 16253.4
$ws.Range("M141").Value = -11073.4

# Sheet2
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 16674742
$ws.Range("I32").Value = 16674742
$ws.Range("K32").Value = 16674742
$ws.Range("M32").Value = -16674455
$ws.Range("H45").Value = 2349.75
$ws.Range("I45").Value = 2349.75
$ws.Range("K45").Value = 2349.75
$ws.Range("M45").Value = -1972.75
$ws.Range("H61").Value = 20883220
$ws.Range("I61").Value = 55560756
$ws.Range("K61").Value = 55560756
$ws.Range("M61").Value = -55560544
$ws.Range("H92").Value = 90539
$ws.Range("J92").Value = 90539
$ws.Range("L92").Value = 90539
$ws.Range("N92").Value = -95531
$ws.Range("H93").Value = 104899
$ws.Range("J93").Value = 104899
$ws.Range("L93").Value = 104899
$ws.Range("N93").Value = -109891
$ws.Range("H101").Value = 109000
$ws.Range("J101").Value = 109000
$ws.Range("L101").Value = 109000
$ws.Range("N101").Value = -115490
$ws.Range("H102").Value = 34797
$ws.Range("I102").Value = 44729.668
$ws.Range("K102").Value = 44729.668
$ws.Range("M102").Value = -43107.668
$ws.Range("H110").Value = 879.5
$ws.Range("I110").Value = 879.5
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 879.5
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1165.5
$ws.Range("N110").Value = ""
$ws.Range("H111").Value = 111995
$ws.Range("J111").Value = 111995
$ws.Range("L111").Value = 111995
$ws.Range("N111").Value = -120175
$ws.Range("H115").Value = 75151
$ws.Range("J115").Value = 75151
$ws.Range("L115").Value = 75151
$ws.Range("N115").Value = -78285
$ws.Range("H122").Value = 1474.75
$ws.Range("I122").Value = 1200
$ws.Range("K122").Value = 3600
$ws.Range("M122").Value = -1150
$ws.Range("H132").Value = 5833.4814
$ws.Range("I132").Value = 2400.0527
$ws.Range("J132").Value = 13987.875
$ws.Range("K132").Value = 7200.158100000001
$ws.Range("L132").Value = 41963.625
$ws.Range("M132").Value = -4670.158100000001
$ws.Range("N132").Value = -47023.625
$ws.Range("H136").Value = 20883220
$ws.Range("I136").Value = 55560756
$ws.Range("K136").Value = 166682268
$ws.Range("M136").Value = -166679718

# Sheet3
$ws = $wb.Worksheets.Item(3)
$ws.Range("H86").Value = 3286.4285
$ws.Range("I86").Value = 3001.25
$ws.Range("J86").Value = 3666.6667
$ws.Range("K86").Value = 3001.25
$ws.Range("L86").Value = 3666.6667
$ws.Range("M86").Value = -1878.25
$ws.Range("N86").Value = -5912.6667
$ws.Range("H89").Value = 3286.4285
$ws.Range("I89").Value = 3001.25
$ws.Range("J89").Value = 3666.6667
$ws.Range("K89").Value = 15006.25
$ws.Range("L89").Value = 18333.3335
$ws.Range("M89").Value = -9390.25
$ws.Range("N89").Value = -29565.3335
$ws.Range("H94").Value = 989.5769
$ws.Range("I94").Value = 989.5769
$ws.Range("K94").Value = 989.5769
$ws.Range("M94").Value = -538.5769
$ws.Range("H134").Value = 31779.676
$ws.Range("I134").Value = 1391.5927
$ws.Range("K134").Value = 4174.7781
$ws.Range("M134").Value = -1639.7781

# Sheet4
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 571245.0600000001
$ws.Range("I31").Value = 1855.5264
$ws.Range("J31").Value = 1062990.5
$ws.Range("K31").Value = 1855.5264
$ws.Range("L31").Value = 1062990.5
$ws.Range("M31").Value = -1560.5264
$ws.Range("N31").Value = -1063580.5
$ws.Range("H34").Value = 571245.0600000001
$ws.Range("I34").Value = 1855.5264
$ws.Range("J34").Value = 1062990.5
$ws.Range("K34").Value = 1855.5264
$ws.Range("L34").Value = 1062990.5
$ws.Range("M34").Value = -1653.5264
$ws.Range("N34").Value = -1063394.5
$ws.Range("H41").Value = 45499
$ws.Range("I41").Value = 15998
$ws.Range("J41").Value = 75000
$ws.Range("K41").Value = 15998
$ws.Range("L41").Value = 75000
$ws.Range("M41").Value = -15570
$ws.Range("N41").Value = -75856
$ws.Range("H50").Value = 60599.8
$ws.Range("J50").Value = 70749.75
$ws.Range("L50").Value = 70749.75
$ws.Range("N50").Value = -71999.75
$ws.Range("H55").Value = 5750
$ws.Range("I55").Value = 5750
$ws.Range("K55").Value = 5750
$ws.Range("M55").Value = -5435
$ws.Range("H58").Value = 2408.7812
$ws.Range("I58").Value = 2266.4348
$ws.Range("K58").Value = 2266.4348
$ws.Range("M58").Value = -2063.4348
$ws.Range("H59").Value = 14000
$ws.Range("I59").Value = 14000
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 14000
$ws.Range("L59").Value = 0
$ws.Range("N59").Value = ""
$ws.Range("M59").Value = -12855
$ws.Range("H68").Value = 66666.664
$ws.Range("J68").Value = 75000
$ws.Range("L68").Value = 75000
$ws.Range("N68").Value = -76498
$ws.Range("H71").Value = 66666.664
$ws.Range("J71").Value = 75000
$ws.Range("L71").Value = 225000
$ws.Range("N71").Value = -232488
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = ""
$ws.Range("N74").Value = ""
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = ""
$ws.Range("N77").Value = ""
$ws.Range("H86").Value = 2829.8572
$ws.Range("I86").Value = 2762
$ws.Range("K86").Value = 2762
$ws.Range("M86").Value = -1639
$ws.Range("H89").Value = 2829.8572
$ws.Range("I89").Value = 2762
$ws.Range("K89").Value = 13810
$ws.Range("M89").Value = -8194
$ws.Range("H104").Value = 62000
$ws.Range("J104").Value = 62000
$ws.Range("L104").Value = 62000
$ws.Range("N104").Value = -67242
$ws.Range("H115").Value = 40909.668
$ws.Range("J115").Value = 40909.668
$ws.Range("L115").Value = 40909.668
$ws.Range("N115").Value = -43259.668
$ws.Range("H116").Value = 54666.5
$ws.Range("J116").Value = 54666.5
$ws.Range("L116").Value = 54666.5
$ws.Range("N116").Value = -63844.5
$ws.Range("H122").Value = 3803.7778
$ws.Range("I122").Value = 3779.25
$ws.Range("K122").Value = 11337.75
$ws.Range("M122").Value = -8887.75
$ws.Range("H132").Value = 2655.4614
$ws.Range("I132").Value = 2365.5454
$ws.Range("J132").Value = 4250
$ws.Range("K132").Value = 7096.6362
$ws.Range("L132").Value = 12750
$ws.Range("M132").Value = -4566.6362
$ws.Range("N132").Value = -17810
$ws.Range("H134").Value = 230259.39
$ws.Range("J134").Value = 5784.5
$ws.Range("L134").Value = 17353.5
$ws.Range("N134").Value = -22423.5
$ws.Range("H136").Value = 2408.7812
$ws.Range("I136").Value = 2266.4348
$ws.Range("K136").Value = 6799.3044
$ws.Range("M136").Value = -4249.3044

# Sheet5
$ws = $wb.Worksheets.Item(5)
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = ""
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").Value = ""
$ws.Range("H54").Value = 5000
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").Value = ""
$ws.Range("H107").Value = 426.83334
$ws.Range("J107").Value = 519.25
$ws.Range("L107").Value = 1557.75
$ws.Range("N107").Value = -5397.75
$ws.Range("H131").Value = 9231.727999999999
$ws.Range("J131").Value = 9079.799999999999
$ws.Range("L131").Value = 27239.4
$ws.Range("N131").Value = -37319.39999999999
$ws.Range("H132").Value = 2073.9546
$ws.Range("I132").Value = 2146.5833
$ws.Range("J132").Value = 1986.8
$ws.Range("K132").Value = 19319.2497
$ws.Range("L132").Value = 17881.2
$ws.Range("M132").Value = -16789.2497
$ws.Range("N132").Value = -22941.2

# Sheet6
$ws = $wb.Worksheets.Item(6)
$ws.Range("H70").Value = 4633.3335
$ws.Range("I70").Value = 4633.3335
$ws.Range("K70").Value = 4633.3335
$ws.Range("M70").Value = -4363.3335
$ws.Range("H73").Value = 4633.3335
$ws.Range("I73").Value = 4633.3335
$ws.Range("K73").Value = 4633.3335
$ws.Range("M73").Value = -3697.3335
$ws.Range("H80").Value = 2513.1667
$ws.Range("I80").Value = 1895
$ws.Range("K80").Value = 1895
$ws.Range("M80").Value = -897
$ws.Range("H83").Value = 2513.1667
$ws.Range("I83").Value = 1895
$ws.Range("K83").Value = 9475
$ws.Range("M83").Value = -4483
$ws.Range("H102").Value = 4364.6665
$ws.Range("I102").Value = 4488.8184
$ws.Range("K102").Value = 4488.8184
$ws.Range("M102").Value = -2866.8184
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = ""
$ws.Range("H113").Value = 4239
$ws.Range("I113").Value = 3948.5715
$ws.Range("J113").Value = 4916.6665
$ws.Range("K113").Value = 3948.5715
$ws.Range("L113").Value = 4916.6665
$ws.Range("M113").Value = -1778.5715
$ws.Range("N113").Value = -9256.666499999999
$ws.Range("H122").Value = 1708.2307
$ws.Range("I122").Value = 1589.7778
$ws.Range("J122").Value = 1974.75
$ws.Range("K122").Value = 4769.3334
$ws.Range("L122").Value = 5924.25
$ws.Range("M122").Value = -2319.3334
$ws.Range("N122").Value = -10824.25
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = ""
$ws.Range("H132").Value = 166670350
$ws.Range("I132").Value = 200003540
$ws.Range("J132").Value = 4398
$ws.Range("K132").Value = 600010620
$ws.Range("L132").Value = 13194
$ws.Range("M132").Value = -600008090
$ws.Range("N132").Value = -18254
$ws.Range("H136").Value = 19607.445
$ws.Range("J136").Value = 19607.445
$ws.Range("L136").Value = 58822.335
$ws.Range("N136").Value = -63922.335

# Sheet7
$ws = $wb.Worksheets.Item(7)
$ws.Range("H40").Value = 3890.5715
$ws.Range("I40").Value = 2719.4443
$ws.Range("J40").Value = 5998.6
$ws.Range("K40").Value = 2719.4443
$ws.Range("L40").Value = 5998.6
$ws.Range("M40").Value = -2583.4443
$ws.Range("N40").Value = -6270.6
$ws.Range("H46").Value = 2920.8076
$ws.Range("I46").Value = 2802.1
$ws.Range("K46").Value = 2802.1
$ws.Range("M46").Value = -2614.1
$ws.Range("H55").Value = 55555956
$ws.Range("I55").Value = 66667092
$ws.Range("J55").Value = 277.66666
$ws.Range("K55").Value = 66667092
$ws.Range("L55").Value = 277.66666
$ws.Range("M55").Value = -66666919
$ws.Range("N55").Value = -623.66666
$ws.Range("H68").Value = 2575.1538
$ws.Range("I68").Value = 2497.7
$ws.Range("J68").Value = 2833.3333
$ws.Range("K68").Value = 2497.7
$ws.Range("L68").Value = 2833.3333
$ws.Range("M68").Value = -1748.7
$ws.Range("N68").Value = -4331.3333
$ws.Range("H71").Value = 2575.1538
$ws.Range("I71").Value = 2497.7
$ws.Range("J71").Value = 2833.3333
$ws.Range("K71").Value = 12488.5
$ws.Range("L71").Value = 14166.6665
$ws.Range("M71").Value = -8744.5
$ws.Range("N71").Value = -21654.6665
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").Value = ""
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").Value = ""
$ws.Range("H93").Value = 71430320
$ws.Range("I93").Value = 76924660
$ws.Range("J93").Value = 4004
$ws.Range("K93").Value = 76924660
$ws.Range("L93").Value = 4004
$ws.Range("M93").Value = -76923412
$ws.Range("N93").Value = -6500
$ws.Range("H103").Value = 35547.25
$ws.Range("J103").Value = 35547.25
$ws.Range("L103").Value = 35547.25
$ws.Range("N103").Value = -37891.25
$ws.Range("H121").Value = 70203.336
$ws.Range("J121").Value = 70203.336
$ws.Range("L121").Value = 70203.336
$ws.Range("N121").Value = -73697.336
$ws.Range("H122").Value = 4572.1763
$ws.Range("I122").Value = 3706
$ws.Range("K122").Value = 11118
$ws.Range("M122").Value = -8668
$ws.Range("H132").Value = 121528.06
$ws.Range("I132").Value = 202222.4
$ws.Range("J132").Value = 87905.414
$ws.Range("K132").Value = 606667.2
$ws.Range("L132").Value = 263716.242
$ws.Range("M132").Value = -604137.2
$ws.Range("N132").Value = -268776.242

# Sheet8
$ws = $wb.Worksheets.Item(8)
$ws.Range("H75").Value = 18813994
$ws.Range("J75").Value = 18813994
$ws.Range("L75").Value = 18813994
$ws.Range("N75").Value = -18815866
$ws.Range("H78").Value = 18813994
$ws.Range("J78").Value = 18813994
$ws.Range("L78").Value = 56441982
$ws.Range("N78").Value = -56451342
$ws.Range("H81").Value = 1199.5
$ws.Range("I81").Value = 1199.5
$ws.Range("K81").Value = 2399
$ws.Range("M81").Value = -1338
$ws.Range("H84").Value = 1199.5
$ws.Range("I84").Value = 1199.5
$ws.Range("K84").Value = 11995
$ws.Range("M84").Value = -6691
$ws.Range("H96").Value = 7666
$ws.Range("J96").Value = 10999
$ws.Range("L96").Value = 10999
$ws.Range("N96").Value = -13745
$ws.Range("H122").Value = 6603.737
$ws.Range("I122").Value = 4798.5557
$ws.Range("K122").Value = 14395.6671
$ws.Range("M122").Value = -11945.6671
$ws.Range("H126").Value = 1162.375
$ws.Range("I126").Value = 1162.375
$ws.Range("K126").Value = 3487.125
$ws.Range("M126").Value = -1017.125
$ws.Range("H130").Value = 88462
$ws.Range("J130").Value = 88462
$ws.Range("L130").Value = 88462
$ws.Range("N130").Value = -98502
$ws.Range("H131").Value = 78832
$ws.Range("J131").Value = 78832
$ws.Range("L131").Value = 78832
$ws.Range("N131").Value = -88912
$ws.Range("H132").Value = 2126.261
$ws.Range("I132").Value = 2145.7
$ws.Range("J132").Value = 1996.6666
$ws.Range("K132").Value = 6437.099999999999
$ws.Range("L132").Value = 5989.9998
$ws.Range("M132").Value = -3907.099999999999
$ws.Range("N132").Value = -11049.9998
$ws.Range("H136").Value = 2191.6667
$ws.Range("I136").Value = 1255.5555
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 3766.6665
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -1216.6665
$ws.Range("N136").Value = -20100
